# controlo de progresso semana 8
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report date header (A2): 2024-12-18 -> 2024-12-19 ---
$ws.Range("A2").Value = 45645

# --- Update "Atribuído a" (assignee) cells with reported progress ---
# E25: Todos -> Francisco Costa
$ws.Range("E25").Value = "Francisco Costa"

# E28: Guilherme Pinho -> Rafael Fernandes
$ws.Range("E28").Value = "Rafael Fernandes"

# E29: Francisco Costa -> Daniel Correia
$ws.Range("E29").Value = "Daniel Correia"

# --- Row 30 (T3.2.3 wrap-up): fill end date, assignee and mark as complete ---
$ws.Range("D30").Value = 45645
$ws.Range("E30").Value = "Rafael Fernandes"
$ws.Range("G30").Value = 1

# --- Row 32 (PT4 / Avaliação): new progress this week ---
# Copy the date formatting used by the neighbouring weekly rows (C27:C29) onto C32
$ws.Range("C27").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C32").Value = 45643
$ws.Range("E32").Value = "Francisco e Guilherme"
$ws.Range("G32").Value = 0.3

# --- Row 33 (T4.1): new progress this week ---
$ws.Range("C27").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C33").Value = 45643
$ws.Range("E33").Value = "Guilherme Pinho"
$ws.Range("G33").Value = 0.2

# --- Row 34 (T4.2): started, date filled in with the same style family as C30 ---
$ws.Range("C30").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C34").Value = 45643

# --- View state: scrolled down a bit and selection moved to J31 ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J31").Select()
